$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2022" column (S) is being appended to the table that previously
# ended at column R. First clone the formatting of the existing last
# column (R, rows 3-14: the header separator row through the bottom
# totals row) into the new column S so number formats / fonts / borders
# stay consistent with the rest of the table, then fill in the 2022
# values.
$ws.Range("R3:R14").Copy() | Out-Null
$ws.Range("S3:S14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Year header
$ws.Range("S4").Value = 2022

# 2022 data values, one per region row
$ws.Range("S5").Value = 27.292394741221504
$ws.Range("S6").Value = 36.613942589338023
$ws.Range("S7").Value = 14.18691257315127
$ws.Range("S8").Value = 55.377118174770182
$ws.Range("S9").Value = 42.247570764681029
$ws.Range("S10").Value = 30.18817294468856
$ws.Range("S11").Value = 97.03085581214826
$ws.Range("S12").Value = 25.2
$ws.Range("S13").Value = 21.849963583394029
$ws.Range("S14").Value = "-"

# Leave the selection the way the author left it after editing
$ws.Range("S16").Select() | Out-Null
